# Update "想去人数" (number of attendees) counts for several events.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) rows
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F6").Value = 544   # 北京·VOCALOID ONLY同人展: 543 -> 544
$wsExpo.Range("F12").Value = 651  # 北京·ICOS SP漫展04动漫节: 649 -> 651
$wsExpo.Range("F19").Value = 159  # 北京·双男主only之皎皎秋月夜: 158 -> 159
$wsExpo.Range("F24").Value = 6562 # 北京·IDO动漫游戏嘉年华47th: 6557 -> 6562
$wsExpo.Range("F25").Value = 4895 # 北京·第19届IJOY漫展xCGF游戏节: 4891 -> 4895
$wsExpo.Range("F35").Value = 605  # 北京·ICOS国际动漫节×CGF中国游戏节04: 603 -> 605

# Sheet "全部类型" (All types) rows - same events, different row positions
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 544   # 北京·VOCALOID ONLY同人展: 543 -> 544
$wsAll.Range("F17").Value = 651   # 北京·ICOS SP漫展04动漫节: 649 -> 651
$wsAll.Range("F24").Value = 159   # 北京·双男主only之皎皎秋月夜: 158 -> 159
$wsAll.Range("F30").Value = 6562  # 北京·IDO动漫游戏嘉年华47th: 6557 -> 6562
$wsAll.Range("F31").Value = 4895  # 北京·第19届IJOY漫展xCGF游戏节: 4891 -> 4895
$wsAll.Range("F38").Value = 605   # 北京·ICOS国际动漫节×CGF中国游戏节04: 603 -> 605
